# Fills out the traceability matrix on "TestCoverageMatrix" with X marks
# for which automated test covers which manual (CMW-T) test case, and
# updates "TestCaseList" accordingly (values only shift because a new
# shared string "X" is introduced).

$wb = $excel.ActiveWorkbook

$matrix = $wb.Worksheets.Item("TestCoverageMatrix")
$list   = $wb.Worksheets.Item("TestCaseList")

# --- TestCoverageMatrix: mark coverage cells with "X" ---
# Columns B..F correspond to CMW-T1..CMW-T5 (row 1 headers).
# Rows 2..7 correspond to AUTO-1, AUTO-3, AUTO-2, AUTO-4, AUTO-6, AUTO-5.
$matrix.Range("B2").Value = "X"   # AUTO-1 covers CMW-T1
$matrix.Range("B3").Value = "X"   # AUTO-3 covers CMW-T1
$matrix.Range("C3").Value = "X"   # AUTO-3 covers CMW-T2
$matrix.Range("C4").Value = "X"   # AUTO-2 covers CMW-T2
$matrix.Range("D5").Value = "X"   # AUTO-4 covers CMW-T3
$matrix.Range("D6").Value = "X"   # AUTO-6 covers CMW-T3
$matrix.Range("F6").Value = "X"   # AUTO-6 covers CMW-T5
$matrix.Range("E7").Value = "X"   # AUTO-5 covers CMW-T4

# Move the active selection like the source workbook ended up with.
$matrix.Activate()
$matrix.Range("O12").Select()
